$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the two new header cells, then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF)
$i0 = @(9, 9, 8, 6, 4, 7, 4, 4, 8, 6, 7)
$if = @(9, 9, 9, 7, 5, 7, 6, 6, 8, 6, 7)

for ($r = 0; $r -lt 11; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
